# Refactor PropertiesExtractor test fixture: reorder the schedule entries
# within each day's block so the subject/day labels line up as expected.
#
# This only changes the displayed text of specific cells; no rows/columns
# are inserted or removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = "Litr1-2"
    "C2" = "Phys1-3"
    "E2" = "Russian1-5"

    "B3" = "Russian1-2"
    "C3" = "Litr1-3"
    "D3" = "Litr1-4"

    "B4" = "Math1-2"
    "D4" = "Math1-4"
    "E4" = "Math1-5"

    "B7" = "Phys2-2"
    "C7" = "Phys2-3"
    "D7" = "English2-4"
    "E7" = "Math2-5"

    "B8" = "Math2-2"
    "C8" = "Math2-3"
    "D8" = "Phys2-4"
    "E8" = "Phys2-5"

    "A11" = "English3-1"
    "B11" = "Phys3-2"
    "C11" = "Russian3-3"
    "D11" = "Litra3-4"
    "E11" = "English3-5"

    "A12" = "Russian3-1"
    "B12" = "Russian3-2"
    "C12" = "Phys3-3"
    "D12" = "Phys3-4"
    "E12" = "Litra3-5"

    "C16" = "Phys4-3"
    "C17" = "Math4-3"

    "D20" = "Phys5-4"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
